# Update trial type data: randomBalloon pump counts (remove unused field in data file)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, "blueBalloon", 119),
    @(3, "greenBalloon", 31),
    @(4, "greenBalloon", 9),
    @(5, "redBalloon", 3),
    @(6, "greenBalloon", 7),
    @(7, "redBalloon", 7),
    @(8, "blueBalloon", 15),
    @(9, "greenBalloon", 19),
    @(10, "greenBalloon", 3),
    @(11, "redBalloon", 6),
    @(12, "blueBalloon", 2),
    @(13, "greenBalloon", 26),
    @(14, "redBalloon", 3),
    @(15, "greenBalloon", 2),
    @(16, "greenBalloon", 1),
    @(17, "greenBalloon", 12),
    @(18, "redBalloon", 2),
    @(19, "redBalloon", 6),
    @(20, "greenBalloon", 19),
    @(21, "redBalloon", 5),
    @(22, "blueBalloon", 108),
    @(23, "blueBalloon", 48),
    @(24, "redBalloon", 5),
    @(25, "blueBalloon", 86),
    @(26, "redBalloon", 4),
    @(27, "blueBalloon", 31),
    @(28, "redBalloon", 3),
    @(29, "blueBalloon", 76),
    @(30, "blueBalloon", 35),
    @(31, "blueBalloon", 69),
    @(32, "redBalloon", 4),
    @(33, "redBalloon", 4),
    @(34, "redBalloon", 4),
    @(35, "redBalloon", 4),
    @(36, "redBalloon", 1),
    @(37, "redBalloon", 5),
    @(38, "redBalloon", 3),
    @(39, "redBalloon", 7),
    @(40, "redBalloon", 7),
    @(41, "redBalloon", 6),
    @(42, "redBalloon", 4),
    @(43, "redBalloon", 1),
    @(44, "redBalloon", 5),
    @(45, "redBalloon", 7),
    @(46, "redBalloon", 5),
    @(47, "redBalloon", 6),
    @(48, "redBalloon", 6),
    @(49, "redBalloon", 5),
    @(50, "redBalloon", 6),
    @(51, "redBalloon", 7),
    @(52, "greenBalloon", 25),
    @(53, "greenBalloon", 6),
    @(54, "greenBalloon", 16),
    @(55, "greenBalloon", 25),
    @(56, "greenBalloon", 27),
    @(57, "greenBalloon", 11),
    @(58, "greenBalloon", 27),
    @(59, "greenBalloon", 10),
    @(60, "greenBalloon", 1),
    @(61, "greenBalloon", 22),
    @(62, "greenBalloon", 8),
    @(63, "greenBalloon", 10),
    @(64, "greenBalloon", 17),
    @(65, "greenBalloon", 9),
    @(66, "greenBalloon", 29),
    @(67, "greenBalloon", 2),
    @(68, "greenBalloon", 2),
    @(69, "greenBalloon", 28),
    @(70, "greenBalloon", 19),
    @(71, "greenBalloon", 13),
    @(72, "blueBalloon", 25),
    @(73, "blueBalloon", 97),
    @(74, "blueBalloon", 73),
    @(75, "blueBalloon", 50),
    @(76, "blueBalloon", 83),
    @(77, "blueBalloon", 125),
    @(78, "blueBalloon", 93),
    @(79, "blueBalloon", 73),
    @(80, "blueBalloon", 5),
    @(81, "blueBalloon", 39),
    @(82, "blueBalloon", 118),
    @(83, "blueBalloon", 53),
    @(84, "blueBalloon", 55),
    @(85, "blueBalloon", 55),
    @(86, "blueBalloon", 27),
    @(87, "blueBalloon", 16),
    @(88, "blueBalloon", 100),
    @(89, "blueBalloon", 86),
    @(90, "blueBalloon", 8),
    @(91, "blueBalloon", 63)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
}
